# This script updates column G ("K") values in the active worksheet,
# replacing the old "Strike#" derived data with regenerated "K" values
# produced by recalculating/re-pulling std/mean and writing s_vals,
# per commit: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2 = 0
    3 = 2
    4 = 0
    5 = 3
    6 = 0
    7 = 2
    8 = 0
    9 = 1
    10 = 2
    11 = 0
    12 = 2
    13 = 2
    14 = 1
    15 = 2
    16 = 2
    17 = 4
    18 = 1
    19 = 3
    20 = 4
    21 = 3
    22 = 0
    23 = 0
    24 = 1
    25 = 2
    26 = 3
    27 = 1
    28 = 3
    29 = 0
    30 = 3
    31 = 2
    32 = 2
    33 = 2
    34 = 2
    35 = 3
    36 = 0
    37 = 0
    38 = 2
    39 = 2
    40 = 3
    41 = 2
    42 = 0
    43 = 2
    44 = 0
    45 = 0
    46 = 1
    47 = 4
    48 = 3
    49 = 2
    50 = 1
    51 = 3
    52 = 5
    53 = 5
    54 = 2
    55 = 3
    56 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
